# cht: add Area chart types to Chart.chart_type
#
# Reposition/resize the three area-chart graphic frames on slide 1 to the
# values captured in the updated fixture. Values are expressed in points
# (PowerPoint COM Left/Top/Width/Height units); the literals below are
# chosen so that, after the engine's internal float32 storage, the
# resulting EMU values land exactly on the target offsets/extents
# (1 pt = 12700 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Chart 1 (standard area chart): off (228600, 228600) ext (2759224, 2768352)
$shp1 = $s.Shapes.Item(1)
$shp1.Left = 18.0
$shp1.Top = 18.0
$shp1.Width = 217.26173400878906
$shp1.Height = 217.98048400878906

# Chart 2 (stacked area chart): off (3200400, 228600) ext (2739752, 2768352)
$shp2 = $s.Shapes.Item(2)
$shp2.Left = 252.0
$shp2.Top = 18.0
$shp2.Width = 215.728515625
$shp2.Height = 217.98048400878906

# Chart 3 (percentStacked area chart): off (6172200, 228600) ext (2720280, 2768352)
$shp3 = $s.Shapes.Item(3)
$shp3.Left = 486.0
$shp3.Top = 18.0
$shp3.Width = 214.19528198242188
$shp3.Height = 217.98048400878906
